# Version 1.0 | Msj: Solucion de multiples if anidados con comprobacion
#
# The vacancy-draw generator used to build this sheet had a bug in its
# nested-IF selection logic. Fixing it (and re-running the draw) produced a
# different set of winning registrations for rows 32-52 while rows 1-31
# (the previously already-confirmed seats) stay untouched.
#
# Apply the corrected draw results directly onto the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# nro_inscripcion, nombre, edad, grado - rows 32..52 (1-based, header in row 1)
$rows = @(
    @(32, 136, "Manuel Vargas",       8,  3),
    @(33, 146, "Victoria Rodríguez",  8,  3),
    @(34, 156, "Manuel Ramos",        8,  3),
    @(35, 164, "Luisa Fernandez",    12,  7),
    @(36, 165, "Manuel Chavez",      10,  5),
    @(37, 168, "Carlos Martinez",     7,  2),
    @(38, 175, "Hugo González",       7,  2),
    @(39, 178, "Julieta Vargas",      8,  3),
    @(40, 179, "Manuel Carril",       7,  2),
    @(41, 181, "Ana Sánchez",        12,  7),
    @(42, 183, "Luisa Torres",       12,  7),
    @(43, 193, "Luisa Peralta",       8,  3),
    @(44, 194, "Luisa Vargas",       12,  7),
    @(45,  74, "Juan Vargas",         8,  3),
    @(46,  97, "Luisa Torres",        8,  3),
    @(47,  98, "Ana Romero",         12,  7),
    @(48, 108, "Victoria Aguilera",   8,  3),
    @(49, 110, "María Fernandez",    12,  7),
    @(50, 114, "Ana Peralta",         8,  3),
    @(51, 126, "Carlos Martinez",    12,  7),
    @(52, 127, "Juan González",       8,  3)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
